$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:R73").AutoFilter() | Out-Null
$name = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$R`$73")
$name.Visible = $False

$ws.Range("A74").Value = "SURVEY_OBSERVATIONS"
$ws.Range("B74").Value = "UOSL"
$ws.Range("C74").Value = "PARITY_checkbox"
$ws.Range("D74").Value = "PARITY"
$ws.Range("E74").Value = "fixed"
$ws.Range("F74").Value = "paritet_5"
$ws.Range("K74").Value = "so_source_column"
$ws.Range("L74").Value = "so_source_value"
$ws.Range("P74").Value = "so_date"
$ws.Range("R74").Value = "keep"
$ws.Range("A75").Value = "SURVEY_OBSERVATIONS"
$ws.Range("B75").Value = "UOSL"
$ws.Range("C75").Value = "obesity_checkbox"
$ws.Range("D75").Value = "OW"
$ws.Range("E75").Value = "fixed"
$ws.Range("F75").Value = "kmi_foer"
$ws.Range("K75").Value = "so_source_column"
$ws.Range("L75").Value = "so_source_value"
$ws.Range("P75").Value = "so_date"
$ws.Range("R75").Value = "keep"
$ws.Range("A76").Value = "SURVEY_OBSERVATIONS"
$ws.Range("B76").Value = "UOSL"
$ws.Range("C76").Value = "blht_checkbox"
$ws.Range("D76").Value = "BLHT"
$ws.Range("E76").Value = "not_fixed"
$ws.Range("F76").Value = "hypertensjon_kronisk"
$ws.Range("G76").Value = 1
$ws.Range("K76").Value = "so_source_column"
$ws.Range("L76").Value = "so_source_value"
$ws.Range("P76").Value = "so_date"
$ws.Range("A77").Value = "SURVEY_OBSERVATIONS"
$ws.Range("B77").Value = "UOSL"
$ws.Range("C77").Value = "bldm_checkbox"
$ws.Range("D77").Value = "BLDM"
$ws.Range("E77").Value = "not_fixed"
$ws.Range("F77").Value = "diabetes_mellitus"
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = 2
$ws.Range("I77").Value = 3
$ws.Range("K77").Value = "so_source_column"
$ws.Range("L77").Value = "so_source_value"
$ws.Range("P77").Value = "so_date"
$ws.Range("A78").Value = "EVENTS"
$ws.Range("B78").Value = "UOSL"
$ws.Range("C78").Value = "obesity_event"
$ws.Range("D78").Value = "BLOB"
$ws.Range("E78").Value = "codesheet"
$ws.Range("F78").Value = "code"
$ws.Range("G78").Value = "vocabulary"
$ws.Range("K78").Value = "event_code"
$ws.Range("L78").Value = "event_record_vocabulary"
$ws.Range("P78").Value = "start_date_record"
$ws.Range("A79").Value = "EVENTS"
$ws.Range("B79").Value = "UOSL"
$ws.Range("C79").Value = "blht_event"
$ws.Range("D79").Value = "BLHT"
$ws.Range("E79").Value = "codesheet"
$ws.Range("F79").Value = "code"
$ws.Range("G79").Value = "vocabulary"
$ws.Range("K79").Value = "event_code"
$ws.Range("L79").Value = "event_record_vocabulary"
$ws.Range("P79").Value = "start_date_record"
$ws.Range("A80").Value = "EVENTS"
$ws.Range("B80").Value = "UOSL"
$ws.Range("C80").Value = "bldm_event"
$ws.Range("D80").Value = "BLDM"
$ws.Range("E80").Value = "codesheet"
$ws.Range("F80").Value = "code"
$ws.Range("G80").Value = "vocabulary"
$ws.Range("K80").Value = "event_code"
$ws.Range("L80").Value = "event_record_vocabulary"
$ws.Range("P80").Value = "start_date_record"
$ws.Range("A81").Value = "EVENTS"
$ws.Range("B81").Value = "UOSL"
$ws.Range("C81").Value = "bldep_event"
$ws.Range("D81").Value = "BLDEP"
$ws.Range("E81").Value = "codesheet"
$ws.Range("F81").Value = "code"
$ws.Range("G81").Value = "vocabulary"
$ws.Range("K81").Value = "event_code"
$ws.Range("L81").Value = "event_record_vocabulary"
$ws.Range("P81").Value = "start_date_record"
$ws.Range("A82").Value = "MEDICINES"
$ws.Range("B82").Value = "UOSL"
$ws.Range("C82").Value = "du_medicines"
$ws.Range("D82").Value = "DU_MED"
$ws.Range("E82").Value = "codesheet"
$ws.Range("F82").Value = "code"
$ws.Range("K82").Value = "medicinal_product_atc_code"
$ws.Range("P82").Value = "date_dispensing"
$ws.Range("A83").Value = "SURVEY_OBSERVATIONS"
$ws.Range("B83").Value = "CHUT"
$ws.Range("C83").Value = "PARITY_checkbox"
$ws.Range("D83").Value = "PARITY"
$ws.Range("E83").Value = "fixed"
$ws.Range("F83").Value = "PARITE"
$ws.Range("K83").Value = "so_source_column"
$ws.Range("L83").Value = "so_source_value"
$ws.Range("P83").Value = "so_date"
$ws.Range("R83").Value = "keep"
$ws.Range("A84").Value = "EVENTS"
$ws.Range("B84").Value = "CHUT"
$ws.Range("C84").Value = "obesity_event"
$ws.Range("D84").Value = "BLOB"
$ws.Range("E84").Value = "codesheet"
$ws.Range("F84").Value = "code"
$ws.Range("G84").Value = "vocabulary"
$ws.Range("K84").Value = "event_code"
$ws.Range("L84").Value = "event_record_vocabulary"
$ws.Range("P84").Value = "start_date_record"
$ws.Range("R84").Value = "keep"
$ws.Range("A85").Value = "EVENTS"
$ws.Range("B85").Value = "CHUT"
$ws.Range("C85").Value = "blht_event"
$ws.Range("D85").Value = "BLHT"
$ws.Range("E85").Value = "codesheet"
$ws.Range("F85").Value = "code"
$ws.Range("G85").Value = "vocabulary"
$ws.Range("K85").Value = "event_code"
$ws.Range("L85").Value = "event_record_vocabulary"
$ws.Range("P85").Value = "start_date_record"
$ws.Range("R85").Value = "keep"
$ws.Range("A86").Value = "EVENTS"
$ws.Range("B86").Value = "CHUT"
$ws.Range("C86").Value = "bldm_event"
$ws.Range("D86").Value = "BLDM"
$ws.Range("E86").Value = "codesheet"
$ws.Range("F86").Value = "code"
$ws.Range("G86").Value = "vocabulary"
$ws.Range("K86").Value = "event_code"
$ws.Range("L86").Value = "event_record_vocabulary"
$ws.Range("P86").Value = "start_date_record"
$ws.Range("R86").Value = "keep"
$ws.Range("A87").Value = "EVENTS"
$ws.Range("B87").Value = "CHUT"
$ws.Range("C87").Value = "bldep_event"
$ws.Range("D87").Value = "BLDEP"
$ws.Range("E87").Value = "codesheet"
$ws.Range("F87").Value = "code"
$ws.Range("G87").Value = "vocabulary"
$ws.Range("K87").Value = "event_code"
$ws.Range("L87").Value = "event_record_vocabulary"
$ws.Range("P87").Value = "start_date_record"
$ws.Range("R87").Value = "keep"
$ws.Range("A88").Value = "MEDICINES"
$ws.Range("B88").Value = "CHUT"
$ws.Range("C88").Value = "du_medicines"
$ws.Range("D88").Value = "DU_MED"
$ws.Range("E88").Value = "codesheet"
$ws.Range("F88").Value = "code"
$ws.Range("K88").Value = "medicinal_product_atc_code"
$ws.Range("P88").Value = "date_dispensing"
$ws.Range("A89").Value = "EVENTS"
$ws.Range("B89").Value = "FISABIO"
$ws.Range("C89").Value = "obesity_event"
$ws.Range("D89").Value = "BLOB"
$ws.Range("E89").Value = "codesheet"
$ws.Range("F89").Value = "code"
$ws.Range("G89").Value = "vocabulary"
$ws.Range("K89").Value = "event_code"
$ws.Range("L89").Value = "event_record_vocabulary"
$ws.Range("P89").Value = "start_date_record"
$ws.Range("A90").Value = "EVENTS"
$ws.Range("B90").Value = "FISABIO"
$ws.Range("C90").Value = "blht_event"
$ws.Range("D90").Value = "BLHT"
$ws.Range("E90").Value = "codesheet"
$ws.Range("F90").Value = "code"
$ws.Range("G90").Value = "vocabulary"
$ws.Range("K90").Value = "event_code"
$ws.Range("L90").Value = "event_record_vocabulary"
$ws.Range("P90").Value = "start_date_record"
$ws.Range("A91").Value = "EVENTS"
$ws.Range("B91").Value = "FISABIO"
$ws.Range("C91").Value = "bldm_event"
$ws.Range("D91").Value = "BLDM"
$ws.Range("E91").Value = "codesheet"
$ws.Range("F91").Value = "code"
$ws.Range("G91").Value = "vocabulary"
$ws.Range("K91").Value = "event_code"
$ws.Range("L91").Value = "event_record_vocabulary"
$ws.Range("P91").Value = "start_date_record"
$ws.Range("A92").Value = "EVENTS"
$ws.Range("B92").Value = "FISABIO"
$ws.Range("C92").Value = "bldep_event"
$ws.Range("D92").Value = "BLDEP"
$ws.Range("E92").Value = "codesheet"
$ws.Range("F92").Value = "code"
$ws.Range("G92").Value = "vocabulary"
$ws.Range("K92").Value = "event_code"
$ws.Range("L92").Value = "event_record_vocabulary"
$ws.Range("P92").Value = "start_date_record"
$ws.Range("A93").Value = "MEDICINES"
$ws.Range("B93").Value = "FISABIO"
$ws.Range("C93").Value = "du_medicines"
$ws.Range("D93").Value = "DU_MED"
$ws.Range("E93").Value = "codesheet"
$ws.Range("F93").Value = "code"
$ws.Range("K93").Value = "medicinal_product_atc_code"
$ws.Range("P93").Value = "date_dispensing"
$ws.Range("A94").Value = "SURVEY_OBSERVATIONS"
$ws.Range("B94").Value = "SAIL"
$ws.Range("C94").Value = "PARITY_checkbox"
$ws.Range("D94").Value = "PARITY"
$ws.Range("E94").Value = "fixed"
$ws.Range("F94").Value = "SERVICE_USER_PARITY_CD"
$ws.Range("K94").Value = "so_source_column"
$ws.Range("L94").Value = "so_source_value"
$ws.Range("P94").Value = "so_date"
$ws.Range("R94").Value = "so_source_value"
$ws.Range("A95").Value = "SURVEY_OBSERVATIONS"
$ws.Range("B95").Value = "SAIL"
$ws.Range("C95").Value = "PBMI_height"
$ws.Range("D95").Value = "PBMI_h"
$ws.Range("E95").Value = "fixed"
$ws.Range("F95").Value = "SERVICE_USER_HEIGHT"
$ws.Range("K95").Value = "so_source_column"
$ws.Range("L95").Value = "so_source_value"
$ws.Range("P95").Value = "so_date"
$ws.Range("R95").Value = "so_source_value"
$ws.Range("A96").Value = "SURVEY_OBSERVATIONS"
$ws.Range("B96").Value = "SAIL"
$ws.Range("C96").Value = "PBMI_weight"
$ws.Range("D96").Value = "PBMI_w"
$ws.Range("E96").Value = "fixed"
$ws.Range("F96").Value = "SERVICE_USER_WEIGHT_KG"
$ws.Range("K96").Value = "so_source_column"
$ws.Range("L96").Value = "so_source_value"
$ws.Range("P96").Value = "so_date"
$ws.Range("R96").Value = "so_source_value"
$ws.Range("A97").Value = "SURVEY_OBSERVATIONS"
$ws.Range("B97").Value = "SAIL"
$ws.Range("C97").Value = "blht_checkbox"
$ws.Range("D97").Value = "BLHT"
$ws.Range("E97").Value = "fixed"
$ws.Range("F97").Value = "EVENT_CD"
$ws.Range("K97").Value = "so_source_column"
$ws.Range("L97").Value = "so_source_value"
$ws.Range("P97").Value = "so_date"
$ws.Range("R97").Value = "so_source_value"
$ws.Range("A98").Value = "SURVEY_OBSERVATIONS"
$ws.Range("B98").Value = "SAIL"
$ws.Range("C98").Value = "bldm_checkbox"
$ws.Range("D98").Value = "BLDM"
$ws.Range("E98").Value = "fixed"
$ws.Range("F98").Value = "EVENT_CD"
$ws.Range("K98").Value = "so_source_column"
$ws.Range("L98").Value = "so_source_value"
$ws.Range("P98").Value = "so_date"
$ws.Range("R98").Value = "so_source_value"
$ws.Range("A99").Value = "EVENTS"
$ws.Range("B99").Value = "SAIL"
$ws.Range("C99").Value = "obesity_event"
$ws.Range("D99").Value = "BLOB"
$ws.Range("E99").Value = "codesheet"
$ws.Range("F99").Value = "code"
$ws.Range("G99").Value = "vocabulary"
$ws.Range("K99").Value = "event_code"
$ws.Range("L99").Value = "event_record_vocabulary"
$ws.Range("P99").Value = "start_date_record"
$ws.Range("R99").Value = "keep"
$ws.Range("A100").Value = "EVENTS"
$ws.Range("B100").Value = "SAIL"
$ws.Range("C100").Value = "blht_event"
$ws.Range("D100").Value = "BLHT"
$ws.Range("E100").Value = "codesheet"
$ws.Range("F100").Value = "code"
$ws.Range("G100").Value = "vocabulary"
$ws.Range("K100").Value = "event_code"
$ws.Range("L100").Value = "event_record_vocabulary"
$ws.Range("P100").Value = "start_date_record"
$ws.Range("R100").Value = "keep"
$ws.Range("A101").Value = "EVENTS"
$ws.Range("B101").Value = "SAIL"
$ws.Range("C101").Value = "bldm_event"
$ws.Range("D101").Value = "BLDM"
$ws.Range("E101").Value = "codesheet"
$ws.Range("F101").Value = "code"
$ws.Range("G101").Value = "vocabulary"
$ws.Range("K101").Value = "event_code"
$ws.Range("L101").Value = "event_record_vocabulary"
$ws.Range("P101").Value = "start_date_record"
$ws.Range("R101").Value = "keep"
$ws.Range("A102").Value = "EVENTS"
$ws.Range("B102").Value = "SAIL"
$ws.Range("C102").Value = "bldep_event"
$ws.Range("D102").Value = "BLDEP"
$ws.Range("E102").Value = "codesheet"
$ws.Range("F102").Value = "code"
$ws.Range("G102").Value = "vocabulary"
$ws.Range("K102").Value = "event_code"
$ws.Range("L102").Value = "event_record_vocabulary"
$ws.Range("P102").Value = "start_date_record"
$ws.Range("R102").Value = "keep"
$ws.Range("A103").Value = "MEDICINES"
$ws.Range("B103").Value = "SAIL"
$ws.Range("C103").Value = "du_medicines"
$ws.Range("D103").Value = "DU_MED"
$ws.Range("E103").Value = "codesheet"
$ws.Range("F103").Value = "code"
$ws.Range("K103").Value = "medicinal_product_atc_code"
$ws.Range("P103").Value = "date_dispensing"
$ws.Range("A104").Value = "SURVEY_OBSERVATIONS"
$ws.Range("B104").Value = "NIHW"
$ws.Range("C104").Value = "PARITY_checkbox"
$ws.Range("D104").Value = "PARITY"
$ws.Range("E104").Value = "fixed"
$ws.Range("F104").Value = "aiemmatsynnytykset"
$ws.Range("K104").Value = "so_source_column"
$ws.Range("L104").Value = "so_source_value"
$ws.Range("P104").Value = "so_date"
$ws.Range("R104").Value = "so_source_value"
$ws.Range("A105").Value = "SURVEY_OBSERVATIONS"
$ws.Range("B105").Value = "NIHW"
$ws.Range("C105").Value = "PBMI_height"
$ws.Range("D105").Value = "PBMI_h"
$ws.Range("E105").Value = "fixed"
$ws.Range("F105").Value = "apituus"
$ws.Range("K105").Value = "so_source_column"
$ws.Range("L105").Value = "so_source_value"
$ws.Range("P105").Value = "so_date"
$ws.Range("R105").Value = "so_source_value"
$ws.Range("A106").Value = "SURVEY_OBSERVATIONS"
$ws.Range("B106").Value = "NIHW"
$ws.Range("C106").Value = "PBMI_weight"
$ws.Range("D106").Value = "PBMI_w"
$ws.Range("E106").Value = "fixed"
$ws.Range("F106").Value = "apaino"
$ws.Range("K106").Value = "so_source_column"
$ws.Range("L106").Value = "so_source_value"
$ws.Range("P106").Value = "so_date"
$ws.Range("R106").Value = "so_source_value"
$ws.Range("A107").Value = "EVENTS"
$ws.Range("B107").Value = "NIHW"
$ws.Range("C107").Value = "obesity_event"
$ws.Range("D107").Value = "BLOB"
$ws.Range("E107").Value = "codesheet"
$ws.Range("F107").Value = "code"
$ws.Range("G107").Value = "vocabulary"
$ws.Range("K107").Value = "event_code"
$ws.Range("L107").Value = "event_record_vocabulary"
$ws.Range("P107").Value = "start_date_record"
$ws.Range("A108").Value = "EVENTS"
$ws.Range("B108").Value = "NIHW"
$ws.Range("C108").Value = "blht_event"
$ws.Range("D108").Value = "BLHT"
$ws.Range("E108").Value = "codesheet"
$ws.Range("F108").Value = "code"
$ws.Range("G108").Value = "vocabulary"
$ws.Range("K108").Value = "event_code"
$ws.Range("L108").Value = "event_record_vocabulary"
$ws.Range("P108").Value = "start_date_record"
$ws.Range("A109").Value = "EVENTS"
$ws.Range("B109").Value = "NIHW"
$ws.Range("C109").Value = "bldm_event"
$ws.Range("D109").Value = "BLDM"
$ws.Range("E109").Value = "codesheet"
$ws.Range("F109").Value = "code"
$ws.Range("G109").Value = "vocabulary"
$ws.Range("K109").Value = "event_code"
$ws.Range("L109").Value = "event_record_vocabulary"
$ws.Range("P109").Value = "start_date_record"
$ws.Range("A110").Value = "EVENTS"
$ws.Range("B110").Value = "NIHW"
$ws.Range("C110").Value = "bldep_event"
$ws.Range("D110").Value = "BLDEP"
$ws.Range("E110").Value = "codesheet"
$ws.Range("F110").Value = "code"
$ws.Range("G110").Value = "vocabulary"
$ws.Range("K110").Value = "event_code"
$ws.Range("L110").Value = "event_record_vocabulary"
$ws.Range("P110").Value = "start_date_record"
$ws.Range("A111").Value = "MEDICINES"
$ws.Range("B111").Value = "NIHW"
$ws.Range("C111").Value = "du_medicines"
$ws.Range("D111").Value = "DU_MED"
$ws.Range("E111").Value = "codesheet"
$ws.Range("F111").Value = "code"
$ws.Range("K111").Value = "medicinal_product_atc_code"
$ws.Range("P111").Value = "date_dispensing"
$ws.Range("A112").Value = "SURVEY_OBSERVATIONS"
$ws.Range("B112").Value = "FERR"
$ws.Range("C112").Value = "PARITY_checkbox"
$ws.Range("D112").Value = "PARITY"
$ws.Range("E112").Value = "fixed"
$ws.Range("F112").Value = "PARTI_PR"
$ws.Range("K112").Value = "so_source_column"
$ws.Range("L112").Value = "so_source_value"
$ws.Range("P112").Value = "so_date"
$ws.Range("R112").Value = "so_source_value"
$ws.Range("A113").Value = "SURVEY_OBSERVATIONS"
$ws.Range("B113").Value = "FERR"
$ws.Range("C113").Value = "PBMI_height"
$ws.Range("D113").Value = "PBMI_h"
$ws.Range("E113").Value = "fixed"
$ws.Range("F113").Value = "ALTEZZA_M"
$ws.Range("K113").Value = "so_source_column"
$ws.Range("L113").Value = "so_source_value"
$ws.Range("P113").Value = "so_date"
$ws.Range("R113").Value = "so_source_value"
$ws.Range("A114").Value = "SURVEY_OBSERVATIONS"
$ws.Range("B114").Value = "FERR"
$ws.Range("C114").Value = "PBMI_weight"
$ws.Range("D114").Value = "PBMI_w"
$ws.Range("E114").Value = "fixed"
$ws.Range("F114").Value = "PESO_M "
$ws.Range("K114").Value = "so_source_column"
$ws.Range("L114").Value = "so_source_value"
$ws.Range("P114").Value = "so_date"
$ws.Range("R114").Value = "so_source_value"
$ws.Range("A115").Value = "EVENTS"
$ws.Range("B115").Value = "FERR"
$ws.Range("C115").Value = "obesity_event"
$ws.Range("D115").Value = "BLOB"
$ws.Range("E115").Value = "codesheet"
$ws.Range("F115").Value = "code"
$ws.Range("G115").Value = "vocabulary"
$ws.Range("K115").Value = "event_code"
$ws.Range("L115").Value = "event_record_vocabulary"
$ws.Range("P115").Value = "start_date_record"
$ws.Range("A116").Value = "EVENTS"
$ws.Range("B116").Value = "FERR"
$ws.Range("C116").Value = "blht_event"
$ws.Range("D116").Value = "BLHT"
$ws.Range("E116").Value = "codesheet"
$ws.Range("F116").Value = "code"
$ws.Range("G116").Value = "vocabulary"
$ws.Range("K116").Value = "event_code"
$ws.Range("L116").Value = "event_record_vocabulary"
$ws.Range("P116").Value = "start_date_record"
$ws.Range("A117").Value = "EVENTS"
$ws.Range("B117").Value = "FERR"
$ws.Range("C117").Value = "bldm_event"
$ws.Range("D117").Value = "BLDM"
$ws.Range("E117").Value = "codesheet"
$ws.Range("F117").Value = "code"
$ws.Range("G117").Value = "vocabulary"
$ws.Range("K117").Value = "event_code"
$ws.Range("L117").Value = "event_record_vocabulary"
$ws.Range("P117").Value = "start_date_record"
$ws.Range("A118").Value = "EVENTS"
$ws.Range("B118").Value = "FERR"
$ws.Range("C118").Value = "bldep_event"
$ws.Range("D118").Value = "BLDEP"
$ws.Range("E118").Value = "codesheet"
$ws.Range("F118").Value = "code"
$ws.Range("G118").Value = "vocabulary"
$ws.Range("K118").Value = "event_code"
$ws.Range("L118").Value = "event_record_vocabulary"
$ws.Range("P118").Value = "start_date_record"
$ws.Range("A119").Value = "MEDICINES"
$ws.Range("B119").Value = "FERR"
$ws.Range("C119").Value = "du_medicines"
$ws.Range("D119").Value = "DU_MED"
$ws.Range("E119").Value = "codesheet"
$ws.Range("F119").Value = "code"
$ws.Range("K119").Value = "medicinal_product_atc_code"
$ws.Range("P119").Value = "date_dispensing"
